# Swap the order of "System" and the other recorder name in the
# "Recorded By" column (column G) of the Session Analysis Results sheet.
# Cells whose value is exactly "System, <other>" (a single other entry)
# become "<other>, System". Cells with additional names (more than one
# comma) or without "System" at all are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    if ($val -match "^System, (.+)$") {
        $rest = $matches[1]
        if ($rest -notmatch ",") {
            $newVal = "$rest, System"
            $cell.Value2 = $newVal
            $changed = $changed + 1
        }
    }
}

Write-Output "Rows changed: $changed"
